# Daily auto-push edit: insert a new reading (2026/02/04, 水, 20, 201)
# right after the existing 2026/02/04 16:00 entry, pushing everything
# from the old row 781 down by one row (old row 822 -> new row 823).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 781; rows 781:822 shift down to 782:823.
$ws.Rows.Item(781).Insert()

# The date (column A) and weekday (column B) are identical to the row
# directly above (780), which is already stored as literal text, not a
# date serial. Copy/PasteSpecial-values from there so the new cells stay
# plain text "2026/02/04" / "水" instead of Excel's COM layer silently
# re-interpreting an assigned "2026/02/04" string as a date value.
$ws.Range("A780:B780").Copy()
$ws.Range("A781:B781").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

# Numeric columns are plain numbers, safe to set directly.
$ws.Cells.Item(781, 3).Value = 20
$ws.Cells.Item(781, 4).Value = 201
